# Update the two-digit-by-two-digit multiplication problems in the
# single table of the document. The table has 20 rows x 5 columns;
# only every 5th row (1, 5, 10, 15, 20) actually holds the problem
# text, the rest are blank "answer" rows.
#
# Because some of the source strings repeat (e.g. "87×52=" appears
# twice but is replaced with two different values), a global
# Find/Replace over the whole story would be ambiguous, so each cell
# is addressed directly by its (row, column) coordinates instead.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; Old = "87×52="; New = "41×52=" },
    @{ Row = 1;  Col = 2; Old = "85×87="; New = "39×47=" },
    @{ Row = 1;  Col = 3; Old = "20×94="; New = "45×62=" },
    @{ Row = 1;  Col = 4; Old = "98×79="; New = "71×32=" },
    @{ Row = 1;  Col = 5; Old = "30×28="; New = "72×16=" },

    @{ Row = 5;  Col = 1; Old = "83×14="; New = "65×60=" },
    @{ Row = 5;  Col = 2; Old = "73×56="; New = "32×54=" },
    @{ Row = 5;  Col = 3; Old = "63×14="; New = "49×52=" },
    @{ Row = 5;  Col = 4; Old = "18×94="; New = "51×16=" },
    @{ Row = 5;  Col = 5; Old = "89×45="; New = "45×95=" },

    @{ Row = 10; Col = 1; Old = "13×17="; New = "41×79=" },
    @{ Row = 10; Col = 2; Old = "36×70="; New = "50×77=" },
    @{ Row = 10; Col = 3; Old = "12×37="; New = "18×42=" },
    @{ Row = 10; Col = 4; Old = "87×52="; New = "41×79=" },
    @{ Row = 10; Col = 5; Old = "31×22="; New = "91×48=" },

    @{ Row = 15; Col = 1; Old = "36×54="; New = "87×95=" },
    @{ Row = 15; Col = 2; Old = "82×80="; New = "35×58=" },
    @{ Row = 15; Col = 3; Old = "73×81="; New = "14×55=" },
    @{ Row = 15; Col = 4; Old = "80×19="; New = "19×54=" },
    @{ Row = 15; Col = 5; Old = "55×98="; New = "83×59=" },

    @{ Row = 20; Col = 1; Old = "50×59="; New = "78×96=" },
    @{ Row = 20; Col = 2; Old = "36×84="; New = "27×69=" },
    @{ Row = 20; Col = 3; Old = "50×11="; New = "26×28=" },
    @{ Row = 20; Col = 4; Old = "93×16="; New = "49×33=" },
    @{ Row = 20; Col = 5; Old = "77×37="; New = "38×33=" }
)

foreach ($edit in $edits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $range = $cell.Range
    # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it
    # before comparing against the plain equation string.
    $current = $range.Text.TrimEnd([char]0x0D, [char]0x07)
    if ($current -ne $edit.Old) {
        throw "Unexpected cell text at row $($edit.Row), col $($edit.Col): expected '$($edit.Old)' but found '$current'"
    }
    $range.Text = $edit.New
}
